# ManageProducts.xlsx - "Manage billing and xml update regarding delivery
# checkbox in manage carrier"
#
# The QA fixture workbook keeps one row per carrier/product on the "Input"
# sheet (rows 2-7, column B = ProductName). Each time the automated test
# suite runs it regenerates a fresh random "prod####" product name for
# every one of those six rows and re-applies the same cell formatting
# (solid white fill + thin top/bottom borders) that the rest of the
# product table already uses.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New product names for this run, keyed by row number.
$newProductNames = @{
    2 = "prodoilI"
    3 = "prodUpeY"
    4 = "prodicfA"
    5 = "prodsfLs"
    6 = "prodgFDI"
    7 = "prodwmQD"
}

foreach ($row in @(2, 3, 4, 5, 6, 7)) {
    $cell = $ws.Range("B" + $row)

    # Write the freshly generated product name for this carrier row.
    $cell.Value = $newProductNames[$row]

    # Re-apply the standard "product name" cell styling: solid white
    # interior fill plus thin borders on the top and bottom edges.
    $cell.Interior.ColorIndex = 2
    $cell.Borders.Item(8).LineStyle = 1
    $cell.Borders.Item(9).LineStyle = 1
}
